# Drop in RMI script results for 3.0
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("SoCDTtiNTY-psgr")

# B2 was a formula (=0.076+(0.076-0.0725)); replace with the plain numeric result
$ws.Range("B2").Value = 0.076

# D2: updated value
$ws.Range("D2").Value = 0.076

# B5 / E5: updated values
$ws.Range("B5").Value = 0.029
$ws.Range("E5").Value = 0.029

# The "About" sheet becomes the active / selected tab
$about = $wb.Worksheets.Item("About")
$about.Activate()
